$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E2E")

# Insert a new row before row 29 (shifts rows 29-62 down to 30-63),
# copying formatting from the row below down into the new row.
[void]$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new SKU test-data entry.
$ws.Range("A29").Value = "SKU-BC128001 -2QTY"
$ws.Range("AI29").Value = "128 oz Oasis"
$ws.Range("AJ29").Value = "review"

# Setting .Value on AJ29 clears the cell's quote-prefix ("Text") format
# that the row-insert had copied down; restore it by pasting the format
# from the neighboring cell that still carries it.
$ws.Range("AJ28").Copy()
$ws.Range("AJ29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# First/last name changes on rows 2 and 3.
$ws.Range("M2").Value = "Hyndavi"
$ws.Range("N2").Value = "Maram"
$ws.Range("M3").Value = "Hyndavi"
$ws.Range("N3").Value = "Maram"

# Update the active selection to match the authored state.
[void]$ws.Range("L11").Select()
